# [Excel Analyzer] Detect named ranges
# Adds a new "Ranges" sheet with an ID/User table and a "Users" named
# range pointing at the User column, mirroring hidden-data patterns like
# external links.

$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet so it lands at the
# end of the tab strip (and becomes the active sheet, like the diff shows).
$sheetCount = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($sheetCount))
$ws.Name = "Ranges"

# Populate the ID / User table data.
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "User"
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Alice"
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Bob"
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Paul"

# Named range over the User column - the kind of hidden-data smell this
# analyzer is meant to flag (alongside external links).
$wb.Names.Add("Users", "=Ranges!`$B`$1:`$B`$4")

# Turn the ID column into a real table named "ID".
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:A4"), $null, 1)
$tbl.Name = "ID"
$tbl.ShowAutoFilter = $false

# The previously active sheet ("External Links") is no longer the
# selected tab now that "Ranges" is active.
$ws.Activate()
